$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 13 (the current "Mean" row), shifting
# the Mean/STD VAR rows down to 15/16.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# The first inserted row (13) stays completely blank/unformatted -- make sure
# it carries no stray content so it collapses out of the saved sheetData.
$ws.Range("A13:C13").ClearContents()

# The second inserted row (14) gets the new "Rising-edge"/"Falling-Edge"
# header labels, formatted like the rest of the bordered table (style used
# by row 12).
$ws.Range("A12:C12").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B14").Value = "Rising-edge"
$ws.Range("C14").Value = "Falling-Edge"

# Update the active selection to match target
$ws.Range("F11").Select()
